# Insert a new weekly price record as row 67, pushing the existing rows
# 67-94 down to 68-95 (dimension grows from A1:R94 to A1:R95).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 67..94 down one position, carrying formatting (e.g. the date
# style on column D) along with them.
$ws.Rows.Item(67).Insert()

# Populate the newly freed row 67 with the new record's data.
$ws.Range("A67").Value = 7
$ws.Range("B67").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C67").Value = "Ñuble"
$ws.Range("D67").Value = 44806
$ws.Range("E67").Value = 16
$ws.Range("F67").Value = 100112031
$ws.Range("G67").Value = "Poroto verde"
$ws.Range("H67").Value = "Magnum"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 60
$ws.Range("K67").Value = 35000
$ws.Range("L67").Value = 35000
$ws.Range("M67").Value = 35000
$ws.Range("N67").Value = "$/malla 25 kilos"
$ws.Range("O67").Value = "Perú"
$ws.Range("P67").Value = 1400
$ws.Range("Q67").Value = 25
$ws.Range("R67").Value = "Hortaliza"
